# issue #5: stock data output to json file
# Insert a new "property_category" column into the "股票" (stock) worksheet,
# with value "stock" for the existing data row, and fix a stray space in the
# company name string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H, shifting old H:J (date, legislator_name,
# legislator_id) one column to the right, to I:K.
$ws.Range("H1:H2").EntireColumn.Insert()

# Fix stray space in the company name text (row2, column B)
$ws.Range("B2").Value = "馬祖酒廠實業股份有限公司(未上市）"

# New header + value for the inserted "property_category" column
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
